$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update Cost ($) and Unit Cost ($/ML) columns ---
$sched = $wb.Worksheets.Item("Schedule")

$sched.Range("E2").Value = 349.94515725
$sched.Range("F2").Value = 7.714840327380953

$sched.Range("E3").Value = 743.421471
$sched.Range("F3").Value = 28.09604954648526

$sched.Range("E4").Value = 13.0052325
$sched.Range("F4").Value = 0.3822819664902998

# --- Sheet "Detailed": update Price (col B) and Type (col C) columns ---
$det = $wb.Worksheets.Item("Detailed")

$det.Range("B30").Value = -5.50985
$det.Range("B31").Value = -2.54301

$det.Range("B32").Value = 0.00003
$det.Range("C32").Value = "historical"

$det.Range("B33").Value = 4.18759
$det.Range("C33").Value = "historical"

$det.Range("B34").Value = 5.18873
$det.Range("B35").Value = -4.98311
$det.Range("B36").Value = 2.2903
$det.Range("B37").Value = 9.87532
$det.Range("B38").Value = 4.54878
$det.Range("B39").Value = 10.64334
$det.Range("B40").Value = 27.59384
$det.Range("B41").Value = 53.90469
$det.Range("B43").Value = 53.90466
$det.Range("B44").Value = 45.19287
$det.Range("B45").Value = 53.99138
$det.Range("B46").Value = 57.01674
$det.Range("B47").Value = 57.03877
$det.Range("B50").Value = 57.06
$det.Range("B52").Value = 56.98
$det.Range("B53").Value = 56.98
$det.Range("B54").Value = 47.77543
$det.Range("B55").Value = 48.13371
$det.Range("B56").Value = 48.7426
$det.Range("B57").Value = 49.58739
$det.Range("B58").Value = 51.23967
$det.Range("B60").Value = 57.06
$det.Range("B61").Value = 58.00516
$det.Range("B62").Value = 60.17548
$det.Range("B64").Value = 35.88
$det.Range("B65").Value = 5.03597
$det.Range("B68").Value = 0.51
$det.Range("B69").Value = 0
$det.Range("B70").Value = 0.51
$det.Range("B71").Value = 0.7
$det.Range("B72").Value = 0.59
$det.Range("B73").Value = 0.51
$det.Range("B74").Value = 0.51
$det.Range("B75").Value = -0.89977
$det.Range("B76").Value = -4.359
$det.Range("B77").Value = -4.20293
$det.Range("B78").Value = -5.58973
$det.Range("B79").Value = -5.66611
$det.Range("B80").Value = -6
$det.Range("B81").Value = -5.58973
$det.Range("B82").Value = -5.43692
$det.Range("B85").Value = -3.21893
$det.Range("B86").Value = 12.21324
$det.Range("B87").Value = 25.73042
$det.Range("B88").Value = 57.06
$det.Range("B92").Value = 55.14746
$det.Range("B93").Value = 56.57874
$det.Range("B94").Value = 56.98
$det.Range("B97").Value = 57.06
